$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Before the edit the table (rows 2-9, columns A-F) looked like:
#   2  median:var
#   3  median:iqr
#   4  median:rvar
#   5  median:skew   <- dropped
#   6  mean:var
#   7  mean:iqr
#   8  mean:rvar
#   9  mean:skew     <- dropped
#
# After the edit the "mean:*" rows come first, the two ":skew" rows (which
# only ever held "nan"/unused skew figures) are removed entirely, and the
# table shrinks to rows 2-7:
#   2  mean:var
#   3  mean:iqr
#   4  mean:rvar
#   5  median:var
#   6  median:iqr
#   7  median:rvar
# ---------------------------------------------------------------------------

# Stash the six rows we keep (everything except the two ":skew" rows) in a
# scratch area below the table first, so that later writes into rows 2-7
# never clobber data that still needs to be read. Copying full A:F ranges
# (rather than whole Rows) keeps the paste scoped to 6 columns and preserves
# both the shared-string typing and the existing cell styles.
$scratchStart = 20
$sourceRows = @(2, 3, 4, 6, 7, 8)
for ($i = 0; $i -lt $sourceRows.Count; $i++) {
    $src = $sourceRows[$i]
    $dst = $scratchStart + $i
    $ws.Range("A$src`:F$src").Copy()
    $ws.Range("A$dst`:F$dst").PasteSpecial(-4163)
    $ws.Range("A$src`:F$src").Copy()
    $ws.Range("A$dst`:F$dst").PasteSpecial(-4122)
}

# Clear out the old table body (rows 2-9) completely.
$ws.Range("A2:F9").Clear()

# Write the scratch rows back in the new order:
#   mean:var, mean:iqr, mean:rvar, median:var, median:iqr, median:rvar
# which, against $sourceRows = (median:var, median:iqr, median:rvar,
# mean:var, mean:iqr, mean:rvar) at scratch offsets (0,1,2,3,4,5), is the
# scratch order (3,4,5,0,1,2).
$scratchOrder = @(3, 4, 5, 0, 1, 2)
for ($i = 0; $i -lt $scratchOrder.Count; $i++) {
    $src = $scratchStart + $scratchOrder[$i]
    $dst = 2 + $i
    $ws.Range("A$src`:F$src").Copy()
    $ws.Range("A$dst`:F$dst").PasteSpecial(-4163)
    $ws.Range("A$src`:F$src").Copy()
    $ws.Range("A$dst`:F$dst").PasteSpecial(-4122)
}

# Remove the scratch rows again.
$ws.Range("A$scratchStart`:F$($scratchStart + $sourceRows.Count - 1)").Clear()

$excel.CutCopyMode = $false
